$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAT_bus")
$ws.Activate()
Write-Output "activated"
